$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 186905.2378535629
$ws.Range("C2").Value = 270692.0082220492
$ws.Range("D2").Value = 316615.9125991673
$ws.Range("E2").Value = 341685.7923241028
$ws.Range("B3").Value = 225589.0428697771
$ws.Range("C3").Value = 324105.1694103569
$ws.Range("D3").Value = 374223.0364166843
$ws.Range("E3").Value = 406075.1521679345
$ws.Range("B4").Value = 201837.4126201456
$ws.Range("C4").Value = 300183.4273028397
$ws.Range("D4").Value = 354764.2568935911
$ws.Range("E4").Value = 387876.1150343845
$ws.Range("B5").Value = 155461.7400494977
$ws.Range("C5").Value = 219824.0359814948
$ws.Range("D5").Value = 246941.4811182787
$ws.Range("E5").Value = 269836.3495078803
$ws.Range("B6").Value = 136799.3467250086
$ws.Range("C6").Value = 192871.4579951396
$ws.Range("D6").Value = 218722.8431517391
$ws.Range("E6").Value = 236088.0884004034
$ws.Range("B7").Value = 14680.21837589119
$ws.Range("C7").Value = 20472.16157971784
$ws.Range("D7").Value = 23277.72962016289
$ws.Range("E7").Value = 24999.89046343153
$ws.Range("B8").Value = 705790.8995429344
$ws.Range("C8").Value = 1030532.188057253
$ws.Range("D8").Value = 1228321.886109981
$ws.Range("E8").Value = 1334288.370835203
$ws.Range("B9").Value = 199574.9481411248
$ws.Range("C9").Value = 280516.8689190554
$ws.Range("D9").Value = 318552.0189318091
$ws.Range("E9").Value = 344320.9075828191
$ws.Range("B10").Value = 85713.91477613166
$ws.Range("C10").Value = 117482.5795151777
$ws.Range("D10").Value = 134894.2441099987
$ws.Range("E10").Value = 143118.8476537207
$ws.Range("B11").Value = 15753.17345700113
$ws.Range("C11").Value = 20501.62131922366
$ws.Range("D11").Value = 23275.28512089467
$ws.Range("E11").Value = 26644.25059606339
$ws.Range("B12").Value = 37561.3566018046
$ws.Range("C12").Value = 53790.5576671302
$ws.Range("D12").Value = 61465.97115054882
$ws.Range("E12").Value = 64367.31937154969
$ws.Range("B13").Value = 48170.31149808713
$ws.Range("C13").Value = 65450.35971744604
$ws.Range("D13").Value = 75873.65920812727
$ws.Range("E13").Value = 81447.00824793526
